# Natmi following Dr Hou advice
# Recompute the Slit2-Robo2 sending/target cluster table across all
# combinations of ECs / FAPs / sCs (3x3 = 9 rows instead of the original 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T
$rows = @(
  ,(2,  @("ECs","Slit2","Robo2","ECs", 3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,3,1,1.655628666666667,4.966886000000001,0.9788840386998182,0.9788840386998183,0.2724905403508889,2.452414863158,0.03412240178955434,0.03412240178955434))
  ,(3,  @("ECs","Slit2","Robo2","FAPs",3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,1,0.3333333333333333,0.01121666666666667,0.03365,0.006631810736596105,0.006631810736596105,0.001846087605555556,0.01661478845,0.0002311747884325316,0.0002311747884325316))
  ,(4,  @("ECs","Slit2","Robo2","sCs", 3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,2,0.6666666666666666,0.02449766666666667,0.073493,0.01448415056358566,0.01448415056358566,0.004031932136555555,0.036287389229,0.0005048953559070445,0.0005048953559070446))
  ,(5,  @("FAPs","Slit2","Robo2","ECs", 3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,3,1,1.655628666666667,4.966886000000001,0.9788840386998182,0.9788840386998183,5.576702051164667,50.190318460482,0.6983378865388631,0.6983378865388633))
  ,(6,  @("FAPs","Slit2","Robo2","FAPs",3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,1,0.3333333333333333,0.01121666666666667,0.03365,0.006631810736596105,0.006631810736596105,0.03778142361666666,0.34003281255,0.004731147419536656,0.004731147419536656))
  ,(7,  @("FAPs","Slit2","Robo2","sCs", 3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,2,0.6666666666666666,0.02449766666666667,0.073493,0.01448415056358566,0.01448415056358566,0.08251620106566666,0.742645809591,0.01033302280249651,0.01033302280249651))
  ,(8,  @("sCs","Slit2","Robo2","ECs", 3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,3,1,1.655628666666667,4.966886000000001,0.9788840386998182,0.9788840386998183,1.967860917532223,17.71074825779,0.2464237503714007,0.2464237503714008))
  ,(9,  @("sCs","Slit2","Robo2","FAPs",3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,1,0.3333333333333333,0.01121666666666667,0.03365,0.006631810736596105,0.006631810736596105,0.01333199913888889,0.11998799225,0.001669488528626917,0.001669488528626917))
  ,(10, @("sCs","Slit2","Robo2","sCs", 3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,2,0.6666666666666666,0.02449766666666667,0.073493,0.01448415056358566,0.01448415056358566,0.02911764079388889,0.262058767145,0.00364623240518211,0.003646232405182112))
)

foreach ($entry in $rows) {
  $r = $entry[0]
  $vals = $entry[1]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
  }
}
